$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = "'26.773.68"
$ws.Range('D2').ClearFormats()

# Row 3
$ws.Range('D3').Value = "'1.644.73"
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  +0.07%  '

# Row 4
$ws.Range('E4').Value = '  +0.58%  '

# Row 5
$ws.Range('D5').Value = "'216.81"
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.47%  '

# Row 6
$ws.Range('D6').Value = "'0.500"
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -0.53%  '

# Row 7
$ws.Range('E7').Value = '  +0.40%  '

# Row 8
$ws.Range('B8').Value = 'Cardano'
$ws.Range('C8').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D8').Value = "'0.251"
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -0.25%  '

# Row 9
$ws.Range('B9').Value = 'Dogecoin'
$ws.Range('C9').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D9').Value = "'0.0629"
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +0.49%  '

# Row 10
$ws.Range('E10').Value = '  -0.25%  '

# Row 11
$ws.Range('E11').Value = '  +0.00%  '

# Row 12
$ws.Range('D12').Value = "'1.867.42"
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -0.26%  '

# Row 13
$ws.Range('D13').Value = "'1.635.91"
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +0.12%  '

# Row 14
$ws.Range('D14').Value = "'4.17"
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -1.23%  '

# Row 15
$ws.Range('E15').Value = '  -0.80%  '

# Row 16
$ws.Range('D16').Value = "'64.48"
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -2.28%  '

# Row 17
$ws.Range('D17').Value = "'26.765.42"
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +0.44%  '

# Row 18
$ws.Range('E18').Value = '  -1.67%  '

# Row 19
$ws.Range('D19').Value = "'213.69"
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -2.19%  '

# Row 20
$ws.Range('E20').Value = '  +0.38%  '

# Row 21
$ws.Range('E21').Value = '  -0.06%  '

# Row 22
$ws.Range('D22').Value = "'2.46"
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +14.35%  '

# Row 24
$ws.Range('E24').Value = '  -2.11%  '

# Row 25
$ws.Range('D25').Value = "'145.11"
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -0.82%  '

# Row 26
$ws.Range('E26').Value = '  +0.43%  '

# Row 27
$ws.Range('E27').Value = '  -1.31%  '

# Row 28
$ws.Range('D28').Value = "'7.12"
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +0.01%  '

# Row 29
$ws.Range('D29').Value = "'15.65"
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -1.33%  '

# Row 30
$ws.Range('E30').Value = '  -1.48%  '

# Row 32
$ws.Range('E32').Value = '  -2.08%  '

# Row 33
$ws.Range('E33').Value = '  -2.06%  '

# Row 34
$ws.Range('D34').Value = "'1.296.03"
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +1.63%  '

# Row 35
$ws.Range('E35').Value = '  -0.21%  '

# Row 36
$ws.Range('E36').Value = '  +1.32%  '

# Row 37
$ws.Range('E37').Value = '  -4.47%  '

# Row 38
$ws.Range('D38').Value = "'0.534"
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +0.86%  '

# Row 39
$ws.Range('D39').Value = "'0.826"
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -0.32%  '

# Row 40
$ws.Range('E40').Value = '  +0.34%  '

# Row 41
$ws.Range('E41').Value = '  +0.04%  '

# Row 42
$ws.Range('E42').Value = '  -0.11%  '

# Row 43
$ws.Range('E43').Value = '  -2.09%  '

# Row 44
$ws.Range('D44').Value = "'1.794.58"
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +0.60%  '

# Row 45
$ws.Range('D45').Value = "'61.90"
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +3.25%  '

# Row 46
$ws.Range('D46').Value = "'91.58"
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -1.68%  '

# Row 47
$ws.Range('E47').Value = '  +0.48%  '

# Row 48
$ws.Range('D48').Value = "'0.0525"
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +1.63%  '

# Row 49
$ws.Range('E49').Value = '  -1.56%  '

# Row 50
$ws.Range('D50').Value = "'0.0976"
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -0.06%  '

# Row 51
$ws.Range('E51').Value = '  +0.23%  '
